$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Barthez"
$ws.Range("B8").Value = "Tonlio"
$ws.Range("C8").Value = "0 201 - 789 - 52 - 30"
$ws.Range("D8").Value = "Ogrenci"
$ws.Range("E8").Value = "ORh+"
$ws.Range("F8").Value = "izmir"
$ws.Range("G8").Value = "buca"
$ws.Range("H8").Value = "Horozluhan Mh izmir / buca`n"
$ws.Range("I8").Value = "müzik dinlemek"

$ws.Rows.Item(8).AutoFit()
